$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new explanation value in E2 (this creates a new shared string entry)
$ws.Range("E2").Value = "Dette er forklaringen"

# Move the active selection to F1 (next empty cell to the right), matching
# the selection change recorded after this edit in the source workbook.
$ws.Range("F1").Select()
